$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-like values (kept as plain strings, no numeric reinterpretation risk)
$ws.Range('D2').Value = '26.840.63'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').Value = '1.875.03'
$ws.Range('E3').Value = '  -1.61%  '
$ws.Range('E5').Value = '  -2.12%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('E9').Value = '  -1.62%  '
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('D13').Value = '1.873.21'
$ws.Range('E13').Value = '  +2.59%  '
$ws.Range('E14').Value = '  -2.72%  '
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '26.867.66'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('E21').Value = '  -2.59%  '
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('E23').Value = '  -1.45%  '
$ws.Range('E24').Value = '  -1.68%  '
$ws.Range('E25').Value = '  -2.65%  '
$ws.Range('E26').Value = '  -0.84%  '
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('E28').Value = '  -2.30%  '
$ws.Range('E29').Value = '  -2.70%  '
$ws.Range('E30').Value = '  -5.35%  '
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('E34').Value = '  -4.50%  '
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('E36').Value = '  +4.88%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E37').Value = '  -3.35%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E38').Value = '  -4.52%  '
$ws.Range('E39').Value = '  -2.26%  '
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('E43').Value = '  -1.63%  '
$ws.Range('E44').Value = '  +3.35%  '
$ws.Range('E45').Value = '  -2.00%  '
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('E47').Value = '  -1.58%  '
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('E49').Value = '  -2.37%  '
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('E51').Value = '  -3.08%  '

# Numeric-looking text values: use a leading apostrophe to force text,
# then reset the style to Normal so no quote-prefix style lingers.
$ws.Range('D4').Formula = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Formula = "'301.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Formula = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Formula = "'0.5366"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Formula = "'0.3752"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Formula = "'0.07192"
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Formula = "'21.60"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Formula = "'0.8904"
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Formula = "'0.08163"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').Formula = "'93.43"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Formula = "'5.323"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Formula = "'14.87"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Formula = "'0.000008538"
$ws.Range('D18').Style = 'Normal'
$ws.Range('D21').Formula = "'4.988"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').Formula = "'6.411"
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Formula = "'2.301"
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Formula = "'146.30"
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Formula = "'18.10"
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Formula = "'1.727"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Formula = "'4.725"
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Formula = "'4.618"
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Formula = "'0.09155"
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Formula = "'0.8128"
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Formula = "'0.05026"
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Formula = "'1.177"
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Formula = "'2.951"
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Formula = "'0.6019"
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Formula = "'2.625"
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Formula = "'3.211"
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Formula = "'0.01955"
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Formula = "'1.070"
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Formula = "'6.631"
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Formula = "'8.894"
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Formula = "'115.22"
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Formula = "'0.5099"
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Formula = "'0.1494"
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Formula = "'10.00"
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Formula = "'1.636"
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Formula = "'37.71"
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Formula = "'0.06055"
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Formula = "'62.22"
$ws.Range('D51').Style = 'Normal'
